# Update the "two-digit divided by one-digit" answer table.
# The table has 20 rows x 5 columns; every 4th row (1,5,9,13,17 in
# 1-based Word indexing) holds the visible equations, the other rows
# are blank spacer rows. We overwrite the 25 equation cells in-place
# (row by row, left to right) with their new values, matching the
# target OOXML exactly while keeping every other part of the document
# (fonts, sizes, paragraph/table formatting, blank rows, etc.)
# untouched.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$newValues = @(
    @("19÷7=2, 5",  "55÷9=6, 1",  "27÷4=6, 3",  "61÷2=30, 1", "31÷6=5, 1"),
    @("67÷9=7, 4",  "77÷2=38, 1", "99÷8=12, 3", "66÷9=7, 3",  "34÷3=11, 1"),
    @("64÷9=7, 1",  "81÷7=11, 4", "40÷6=6, 4",  "28÷5=5, 3",  "57÷6=9, 3"),
    @("26÷8=3, 2",  "88÷9=9, 7",  "80÷6=13, 2", "46÷7=6, 4",  "96÷8=12, 0"),
    @("84÷3=28, 0", "11÷4=2, 3",  "14÷7=2, 0",  "44÷8=5, 4",  "90÷7=12, 6")
)

$dataRows = @(1, 5, 9, 13, 17)

for ($i = 0; $i -lt $dataRows.Length; $i++) {
    $row = $dataRows[$i]
    $values = $newValues[$i]
    for ($c = 1; $c -le 5; $c++) {
        $cell = $t.Cell($row, $c)
        $cell.Range.Text = $values[$c - 1]
    }
}
